# interstellar.xlsx maintenance pass:
#  - Routes/Traffic: planet "H" was a typo/placeholder that doesn't exist on the
#    "Planet Names" sheet. Route 4 (B->H) should read B->E, and Route 5's origin
#    should read C instead of B.
#  - Selections/active tab move to reflect where the editor was working (Traffic
#    tab ends up active/selected, with Routes' and Traffic's last selections
#    updated too).

$wb = $excel.ActiveWorkbook

$routes  = $wb.Worksheets.Item("Routes")
$traffic = $wb.Worksheets.Item("Traffic")

# --- Fix the bad "H" planet reference / route origin on both sheets ---
$routes.Range("C5").Value = "E"
$routes.Range("B6").Value = "C"

$traffic.Range("C5").Value = "E"
$traffic.Range("B6").Value = "C"

# --- Update selections (Routes first, Traffic last so Traffic ends up active) ---
$routes.Range("C6").Select() | Out-Null

$traffic.Activate() | Out-Null
$traffic.Range("B6").Select() | Out-Null
